# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values on the active sheet to reflect
# the regenerated data, leaving all other columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 3
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 3
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 4
    21 = 2
    22 = 3
    24 = 3
    25 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
